$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at FX (pushing the old FX:GV block — punt_returns..EXP —
# one column to the right, to FY:GW), then label the new column's header.
$ws.Columns("FX:FX").Insert()
$ws.Range("FX1").Value = "hit_within3years"

# Match the saved selection state (active cell moved to FW3).
$ws.Range("FW3").Select() | Out-Null
